# Update the "datetimeFigureOut" Date placeholder text from 12/11/21 to
# 12/13/21 everywhere it appears: on the slide master and on every one of
# its slide layouts (PowerPoint stamps this placeholder's cached text on
# the master + each layout independently).

$p = $ppt.ActivePresentation
$newDate = "12/13/21"

# -- Slide Master --------------------------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# -- Every Slide Layout off the master ------------------------------------
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}
